# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the four "low" priority
# rows (0023b088, 179e0ac6, 33a44c41, 89e4f584) in both the zh-cn and
# de-de sheets are re-classified from Priority "low" to "ht", and their
# "Latest Handoff Datetime" is bumped forward a few seconds to reflect
# the freshly generated handoff xliff.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: rows 4-7 are the "low" priority files.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
for ($r = 4; $r -le 7; $r++) {
    $wsZhCn.Cells.Item($r, 5).Value = "ht"                     # E: Priority
    $wsZhCn.Cells.Item($r, 8).Value = "2016-09-05 10:43:29"    # H: Latest Handoff Datetime
}

# de-de sheet: rows 4-7 are the "low" priority files.
$wsDeDe = $wb.Worksheets.Item("de-de")
for ($r = 4; $r -le 7; $r++) {
    $wsDeDe.Cells.Item($r, 5).Value = "ht"                     # E: Priority
    $wsDeDe.Cells.Item($r, 8).Value = "2016-09-05 10:43:34"    # H: Latest Handoff Datetime
}

# Overview sheet mirrors de-de's "Latest Handoff Datetime" in its
# "Latest HO Xliff Generate Date" column, so it picks up the same bump.
$wsOverview = $wb.Worksheets.Item("Overview")
for ($r = 4; $r -le 7; $r++) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-09-05 10:43:34"   # G: Latest HO Xliff Generate Date
}
